$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.175700306892395
$ws.Range("B1").Value = 2.405180692672729
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.347479581832886
$ws.Range("E1").Value = 1.207398533821106
